# Add entries for the last three days (28, 29, 30 Nov 2020) of the
# NOV-2020 daily tracker sheet, matching the formatting used by the
# preceding rows, and update the sheet's saved selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NOV-2020")
$wsOct = $wb.Worksheets.Item("OCT-2020")

$fmtFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---- Row 29 (28 Nov 2020) - Completed ----
$ws.Range("C28:F28").Copy() | Out-Null
$ws.Range("C29:F29").PasteSpecial($fmtFormats) | Out-Null
$ws.Range("C29").Value = "Sonia Application"
$ws.Range("D29").Value = "Regression testing on Sonia Application(Best ivc report)"
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = "Completed"
$ws.Rows.Item(29).RowHeight = 30

# ---- Row 30 (29 Nov 2020) - Completed ----
$ws.Range("C28:F28").Copy() | Out-Null
$ws.Range("C30:F30").PasteSpecial($fmtFormats) | Out-Null
$ws.Range("C30").Value = "Sonia Application"
$ws.Range("D30").Value = "Regression testing on Sonia Application(Best ivc report)"
$ws.Range("E30").Value = 1
$ws.Range("F30").Value = "Completed"
$ws.Rows.Item(30).RowHeight = 30

# ---- Row 31 (30 Nov 2020) - WIP ----
$ws.Range("C28:D28").Copy() | Out-Null
$ws.Range("C31:D31").PasteSpecial($fmtFormats) | Out-Null
$ws.Range("E28").Copy() | Out-Null
$ws.Range("E31").PasteSpecial($fmtFormats) | Out-Null
$wsOct.Range("F2").Copy() | Out-Null
$ws.Range("F31").PasteSpecial($fmtFormats) | Out-Null
$ws.Range("C31").Value = "Sonia Application"
$ws.Range("D31").Value = "Sanity testing on B2C app, QMVAR site, GSS site and Hayaai site. Regression testing on Sonia Application(All Denka, Best Denki, Bic Camera, Homac, Koroganeya, MrMax and Stream)"
$ws.Range("E31").ClearContents() | Out-Null
$ws.Range("F31").Value = "WIP"
$ws.Rows.Item(31).RowHeight = 30

# ---- Update the sheet's saved selection ----
$ws.Activate() | Out-Null
$ws.Range("D34").Select() | Out-Null
